# 1st changes of mifos to finflux
#
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet - this pushes the old "Late"/"Waived"/"Outstanding" columns one
# slot to the right (N->O, O->P, P->Q) and leaves the new column N empty.
# Also makes "Repayment schedule" the active sheet/tab and updates its
# selection, matching the author's last interaction before saving.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column in front of column N (shifts N:P -> O:Q).
$wsSchedule.Columns("N").Insert()

# Make "Repayment schedule" the active sheet/tab (was "Edit Repayment
# Schedule" before) and move the selection to L13.
$wsSchedule.Activate()
$wsSchedule.Range("L13").Select()
